# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right before the existing "2022-Q3"
#    sheet (so the sheet order becomes 总计, 2022-Q4, 2022-Q3, 2022-Q2) and
#    populate it with the Q4 fund-holding data.
# 2. Insert a new summary row for "2022-Q4" at the top of the data in the
#    "总计" sheet, shifting the existing 2022-Q3 / 2022-Q2 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q4" worksheet, positioned before "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row (bold/centered/bordered style copied from the 2022-Q3 sheet)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Copy the header style (B1:H1) from the 2022-Q3 sheet so formatting matches.
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# Row 2 data
$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "006729"
$q4.Range("C2").NumberFormat = "@"
$q4.Range("C2").Value = "万家中证500指数增强A"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "3.13"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "93.56"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "1.32"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0413"
$q4.Range("H2").Value = 3

# Row 3 data
$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "006730"
$q4.Range("C3").NumberFormat = "@"
$q4.Range("C3").Value = "万家中证500指数增强C"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "2.38"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "93.56"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "1.32"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0314"
$q4.Range("H3").Value = 3

# Copy the A2 style (used for the row-index column) onto A2:A3 of the new sheet.
$q3.Range("A2").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)

$q4.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 2: insert the "2022-Q4" summary row into the "总计" sheet
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Shift the existing data (old row 2 -> row 3, old row 3 -> row 4) down,
# then overwrite row 2 with the new 2022-Q4 summary values.
$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2022-Q2"
$zj.Range("C4").Value = 2
$zj.Range("D4").Value = 0.99

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2022-Q3"
$zj.Range("C3").Value = 1
$zj.Range("D3").Value = 0.01

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 0.07000000000000001

# Ensure the row-index column keeps its original style on the newly
# extended row 4.
$zj.Range("A2").Copy()
$zj.Range("A4").PasteSpecial(-4122)

$zj.Range("A1").Select()
